# Add three new rows to the "file format" master data table for the
# new "html" file format (English, Arabic, French descriptions).
#
# New rows are entered in this order (Arabic, then English, then French)
# so that the new shared-string entries end up with the same index
# assignment as in the target workbook:
#   23 = "html"
#   24 = "ملف html"
#   25 = "html file"
#   26 = "Fichier html"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$TRUE_ALIGN_LEFT = -4131   # xlLeft

# Row 12 - Arabic "html" entry (filled first)
$ws.Range("A12").Value = "html"
$ws.Range("B12").Value = "ملف html"
$ws.Range("C12").Value = "ara"
$ws.Range("D12").Value = $true
$ws.Range("D12").HorizontalAlignment = $TRUE_ALIGN_LEFT
$ws.Range("E12").Value = "superadmin"
$ws.Range("F12").Value = "now()"

# Row 11 - English "html" entry (filled second)
$ws.Range("A11").Value = "html"
$ws.Range("B11").Value = "html file"
$ws.Range("C11").Value = "eng"
$ws.Range("D11").Value = $true
$ws.Range("D11").HorizontalAlignment = $TRUE_ALIGN_LEFT
$ws.Range("E11").Value = "superadmin"
$ws.Range("F11").Value = "now()"

# Row 13 - French "html" entry (filled third)
$ws.Range("A13").Value = "html"
$ws.Range("B13").Value = "Fichier html"
$ws.Range("C13").Value = "fra"
$ws.Range("D13").Value = $true
$ws.Range("D13").HorizontalAlignment = $TRUE_ALIGN_LEFT
$ws.Range("E13").Value = "superadmin"
$ws.Range("F13").Value = "now()"

# Leave the selection where it ends up after working in columns G onward,
# matching the saved view state of the target workbook.
$ws.Range("G1:XFD1048576").Select()
